$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.411.30"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.683.04"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "683.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").Value = "3.681.86"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.58%  "
$ws.Range("D14").Value = "4.303.79"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.01%  "
$ws.Range("D16").Value = "3.682.84"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "69.388.43"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.28%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.66%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.646"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("D25").Value = "3.830.13"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.67%  "
$ws.Range("E32").Value = "  -9.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.23%  "
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").Value = "3.657.27"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0898"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.83%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.942"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.92%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "165.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -12.45%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000275"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "28.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.29%  "
